$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$np = $s.NotesPage
$sh2 = $np.Shapes.Item(2)
$sh3 = $np.Shapes.Item(3)
Write-Output ("sh2 W=" + $sh2.Width + " H=" + $sh2.Height)
Write-Output ("sh3 W=" + $sh3.Width + " H=" + $sh3.Height)
